$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every date in column A (rows 1-366) forward by 8 days.
$start = Get-Date -Year 2025 -Month 3 -Day 2
for ($row = 1; $row -le 366; $row++) {
    $d = $start.AddDays($row - 1)
    $s = $d.ToString("yyyy-MM-dd")
    # Leading apostrophe forces the value to be stored as text instead of
    # being auto-converted to a date serial number by Excel.
    $ws.Cells.Item($row, 1).Value = "'" + $s
}

# Update the notes for the row now showing 2025-03-03.
$ws.Range("B2").Value = "fawef"
$ws.Range("C2").Value = "wef"
